$d = $word.ActiveDocument

# Locate the very end of the document body: the end of the final
# "git push" paragraph that is already present in the document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

# Raw WordprocessingML for the paragraphs being appended, exactly as
# authored: two blank lines, "py generate_index.py", "git add .",
# the "git commit" line, "git pull --no-rebase", "git push", and a
# trailing blank paragraph -- including the proofing marks and the
# east-Asia font hints Word records for text typed this way.
$newParagraphsXml = '<w:p/><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>py</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> generate_index.py</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">git </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>add .</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>git commit -m "</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>update grouped full index</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>"</w:t></w:r></w:p><w:p><w:r><w:t>git pull --no-rebase</w:t></w:r></w:p><w:p><w:r><w:t>git push</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'

$wordXmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$wordXmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$wordXml = $wordXmlHeader + $newParagraphsXml + $wordXmlFooter

$insertionPoint.InsertXML($wordXml)
